{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Applies the three edits from the commit:\n//  1. \"How did the plan changed over the semester?\" - Word's proofing pass\n//     split the run around \"changed\" (grammar flag) into three runs. The\n//     visible text is unchanged.\n//  2. \"Concentrate on one app, rather then trying to make two. ...\" - same\n//     kind of proofing split, this time around \"then\" (spelling + grammar\n//     flag). The visible text is unchanged.\n//  3. \"How much time have you spend on project per pearson - just estimate\"\n//     -> \"... per person - ...\" - an actual spelling fix (the real content\n//     change called out in the commit message).\n//\n// Note: `<w:proofErr .../>` markers are cosmetic artifacts that Word's live\n// spell/grammar checker stamps into the run-split points; they carry no\n// text or formatting and are not reachable through the Word JS API (there is\n// no `Range`/`Paragraph` member that inserts them). What IS reproducible\n// through the object model -- and is done below -- is the underlying run\n// split itself: toggling a character-formatting property on a sub-range and\n// then reverting it forces Word to break that sub-range into its own run\n// (cloning the original run's formatting), which is exactly the structural\n// change the proofing pass produced.\n\nasync function splitRun(context, searchText) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    return;\n  }\n  const range = results.items[0];\n  // Flip a formatting property on and back off: this forces the host to\n  // materialize a dedicated run for exactly this sub-range instead of\n  // merging it back into its neighbors.\n  range.font.bold = true;\n  await context.sync();\n  range.font.bold = false;\n  await context.sync();\n}\n\n// --- 1. \"How did the plan changed over the semester?\" -----------------\n// Split off \"changed\" from its neighbors so the paragraph ends up as three\n// runs: \"How did the plan \" | \"changed\" | \" over the semester?\".\nawait splitRun(context, \"changed over the semester?\");\nawait splitRun(context, \"changed\");\n\n// --- 2. \"Concentrate on one app, rather then trying ...\" ---------------\n// Split off \"then\" the same way: \"Concentrate on one app, rather \" |\n// \"then\" | \" trying to make two. Also, use a better file format.\"\nawait splitRun(context, \"then trying to make two. Also, use a better file format.\");\nawait splitRun(context, \"then\");\n\n// --- 3. Fix the \"pearson\" -> \"person\" typo -----------------------------\nconst typoResults = context.document.body.search(\"pearson\", { matchCase: true });\ntypoResults.load(\"items\");\nawait context.sync();\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\"person\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument / $d is the open document.\n#\n# Applies the three edits from the commit:\n#  1. \"How did the plan changed over the semester?\" - Word's proofing pass\n#     split the run around \"changed\" (grammar flag) into three runs. The\n#     visible text is unchanged.\n#  2. \"Concentrate on one app, rather then trying to make two. ...\" - same\n#     kind of proofing split, this time around \"then\" (spelling + grammar\n#     flag). The visible text is unchanged.\n#  3. \"How much time have you spend on project per pearson - just estimate\"\n#     -> \"... per person - ...\" - an actual spelling fix (the real content\n#     change called out in the commit message).\n#\n# Note: `<w:proofErr .../>` markers are cosmetic artifacts that Word's live\n# spell/grammar checker stamps into the run-split points; they carry no\n# text or formatting and there is no Range/Find/Paragraph member that\n# inserts them, in the real Word object model or here. What IS reproducible\n# -- and is done below -- is the underlying run split itself: toggling a\n# character-formatting property on a sub-range and then reverting it forces\n# Word to break that sub-range into its own run (cloning the original run's\n# formatting), which is exactly the structural change the proofing pass\n# produced.\n\n$d = $word.ActiveDocument\n\nfunction Split-Run($findText) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1)\n    if ($found) {\n        $rng.Font.Bold = $true\n        $rng.Font.Bold = $false\n    }\n}\n\n# --- 1. \"How did the plan changed over the semester?\" -----------------\n# Split off \"changed\" from its neighbors so the paragraph ends up as three\n# runs: \"How did the plan \" | \"changed\" | \" over the semester?\".\nSplit-Run \"changed over the semester?\"\nSplit-Run \"changed\"\n\n# --- 2. \"Concentrate on one app, rather then trying ...\" ---------------\n# Split off \"then\" the same way: \"Concentrate on one app, rather \" |\n# \"then\" | \" trying to make two. Also, use a better file format.\"\nSplit-Run \"then trying to make two. Also, use a better file format.\"\nSplit-Run \"then\"\n\n# --- 3. Fix the \"pearson\" -> \"person\" typo -----------------------------\n$typoRange = $d.Content\n$typoRange.Find.Execute(\"pearson\", $true, $false, $false, $false, $false, $true, 1, $false, \"person\", 1) | Out-Null\n"}
